$d = $word.ActiveDocument

$newShapesXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251667456" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3A370A83" wp14:editId="35150568"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>1809750</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>4467225</wp:posOffset></wp:positionV><wp:extent cx="3600450" cy="3952875"/><wp:effectExtent l="0" t="0" r="19050" b="28575"/><wp:wrapNone/><wp:docPr id="5" name="Rectangle 5"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="3600450" cy="3952875"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FFBA00"/></a:solidFill></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect w14:anchorId="4F2A9304" id="Rectangle 5" o:spid="_x0000_s1026" style="position:absolute;margin-left:142.5pt;margin-top:351.75pt;width:283.5pt;height:311.25pt;z-index:251667456;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQDEiQLLjwIAAG4FAAAOAAAAZHJzL2Uyb0RvYy54bWysVEtv2zAMvg/YfxB0X+2kSR9BnSJrkWFA&#10;0RZth54VWYoNyKJGKXGyXz9KdtygK3YYloNCmuTHN6+ud41hW4W+Blvw0UnOmbISytquC/7jZfnl&#10;gjMfhC2FAasKvleeX88/f7pq3UyNoQJTKmQEYv2sdQWvQnCzLPOyUo3wJ+CUJaEGbEQgFtdZiaIl&#10;9MZk4zw/y1rA0iFI5T19ve2EfJ7wtVYyPGjtVWCm4BRbSC+mdxXfbH4lZmsUrqplH4b4hygaUVty&#10;OkDdiiDYBus/oJpaInjQ4URCk4HWtVQpB8pmlL/L5rkSTqVcqDjeDWXy/w9W3m8fkdVlwaecWdFQ&#10;i56oaMKujWLTWJ7W+RlpPbtH7DlPZMx1p7GJ/5QF26WS7oeSql1gkj6enuX5ZEqVlyQ7vZyOL84T&#10;avZm7tCHbwoaFomCI7lPpRTbOx/IJakeVKI3D6Yul7UxicH16sYg2wrq73L5dZGnlpLJkVoWU+iC&#10;TlTYGxWNjX1SmnKnMMfJY5o6NeAJKZUNo05UiVJ1bqY5/WJlopc4p9EicQkwImsKb8DuAQ6aHcgB&#10;u4Pp9aOpSkM7GOd/C6wzHiySZ7BhMG5qC/gRgKGses+dPoV/VJpIrqDc02QgdCvjnVzW1J874cOj&#10;QNoR6intfXigRxtoCw49xVkF+Ouj71GfRpeknLW0cwX3PzcCFWfmu6WhvhxNJnFJEzOZno+JwWPJ&#10;6lhiN80NUNtHdGGcTGTUD+ZAaoTmlc7DInolkbCSfBdcBjwwN6G7BXRgpFoskhotphPhzj47GcFj&#10;VeP8vexeBbp+SAPN9z0c9lPM3s1qpxstLSw2AXSdBvmtrn29aanT4PQHKF6NYz5pvZ3J+W8AAAD/&#10;/wMAUEsDBBQABgAIAAAAIQCyl0Cm4QAAAAwBAAAPAAAAZHJzL2Rvd25yZXYueG1sTI/LTsMwEEX3&#10;SPyDNUjsqE2qFCvEqRAVUlZIpCxg58ZDEtWPKHbbhK9nWMFyZo7unFtuZ2fZGac4BK/gfiWAoW+D&#10;GXyn4H3/cieBxaS90TZ4VLBghG11fVXqwoSLf8NzkzpGIT4WWkGf0lhwHtsenY6rMKKn21eYnE40&#10;Th03k75QuLM8E2LDnR48fej1iM89tsfm5BTsrK139VHq72X5dEJ+YNbUr0rd3sxPj8ASzukPhl99&#10;UoeKnA7h5E1kVkEmc+qSFDyIdQ6MCJlntDkQus42AnhV8v8lqh8AAAD//wMAUEsBAi0AFAAGAAgA&#10;AAAhALaDOJL+AAAA4QEAABMAAAAAAAAAAAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwEC&#10;LQAUAAYACAAAACEAOP0h/9YAAACUAQAACwAAAAAAAAAAAAAAAAAvAQAAX3JlbHMvLnJlbHNQSwEC&#10;LQAUAAYACAAAACEAxIkCy48CAABuBQAADgAAAAAAAAAAAAAAAAAuAgAAZHJzL2Uyb0RvYy54bWxQ&#10;SwECLQAUAAYACAAAACEAspdApuEAAAAMAQAADwAAAAAAAAAAAAAAAADpBAAAZHJzL2Rvd25yZXYu&#10;eG1sUEsFBgAAAAAEAAQA8wAAAPcFAAAAAA==&#10;" fillcolor="#ffba00" strokecolor="#1f3763 [1604]" strokeweight="1pt"/></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251665408" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="1FC5F68E" wp14:editId="027D2D71"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>-333375</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>5095875</wp:posOffset></wp:positionV><wp:extent cx="3600450" cy="3952875"/><wp:effectExtent l="0" t="0" r="19050" b="28575"/><wp:wrapNone/><wp:docPr id="4" name="Rectangle 4"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="3600450" cy="3952875"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FFD869"/></a:solidFill></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect w14:anchorId="1C8F523D" id="Rectangle 4" o:spid="_x0000_s1026" style="position:absolute;margin-left:-26.25pt;margin-top:401.25pt;width:283.5pt;height:311.25pt;z-index:251665408;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQD1kLAMjQIAAG4FAAAOAAAAZHJzL2Uyb0RvYy54bWysVEtv2zAMvg/YfxB0X+28+gjqFEGLDgOK&#10;tmg79KzIUmJAFjVKiZP9+lGy4wZtscOwi0ya5Mc3L692jWFbhb4GW/DRSc6ZshLK2q4L/vPl9ts5&#10;Zz4IWwoDVhV8rzy/Wnz9ctm6uRpDBaZUyAjE+nnrCl6F4OZZ5mWlGuFPwClLQg3YiEAsrrMSRUvo&#10;jcnGeX6atYClQ5DKe/p70wn5IuFrrWR40NqrwEzBKbaQXkzvKr7Z4lLM1yhcVcs+DPEPUTSituR0&#10;gLoRQbAN1h+gmloieNDhREKTgda1VCkHymaUv8vmuRJOpVyoON4NZfL/D1bebx+R1WXBJ5xZ0VCL&#10;nqhowq6NYpNYntb5OWk9u0fsOU9kzHWnsYlfyoLtUkn3Q0nVLjBJPyeneT6dUeUlySYXs/H52Syi&#10;Zm/mDn34rqBhkSg4kvtUSrG986FTPahEbx5MXd7WxiQG16trg2wrqL835CtPLSX0I7UsptAFnaiw&#10;NyoaG/ukNOVOYY6TxzR1asATUiobRp2oEqXq3MyOvcQ5jRYpowQYkTWFN2D3AAfNDuSA3eXX60dT&#10;lYZ2MM7/FlhnPFgkz2DDYNzUFvAzAENZ9Z47fQr/qDSRXEG5p8lA6FbGO3lbU3/uhA+PAmlHqKe0&#10;9+GBHm2gLTj0FGcV4O/P/kd9Gl2SctbSzhXc/9oIVJyZH5aG+mI0ncYlTcx0djYmBo8lq2OJ3TTX&#10;QG0f0YVxMpFRP5gDqRGaVzoPy+iVRMJK8l1wGfDAXIfuFtCBkWq5TGq0mE6EO/vsZASPVY3z97J7&#10;Fej6IQ003/dw2E8xfzernW60tLDcBNB1GuS3uvb1pqVOg9MfoHg1jvmk9XYmF38AAAD//wMAUEsD&#10;BBQABgAIAAAAIQCOIdNB3wAAAAwBAAAPAAAAZHJzL2Rvd25yZXYueG1sTI/BTsMwDIbvk3iHyEjc&#10;tqTb2pXSdEJIIDjssDGBuGVNaCsSp2qyrvD0mBMcbX/6/f3ldnKWjWYIvUcJyVwAM9h43WMr4fXl&#10;YVYAC1GhVtajkfBlAmyrq1mpCu0vuDfjIbaMQjAUSkIX41hwHprOOBXmfjRIt48/OBVpHFquB3Wh&#10;cGf5QoiMO9UjfejUaO4703weTk7C2PV5vrNPdfbm3r63q8fnQaTvUt5cT3e3wKKZ4h8Mv/qkDhU5&#10;Hf0JdWBWwizNiJQwS9dLYARkyYaWI6GLfCOAVyX/36H6AQAA//8DAFBLAQItABQABgAIAAAAIQC2&#10;gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAG&#10;AAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAG&#10;AAgAAAAhAPWQsAyNAgAAbgUAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1sUEsBAi0A&#10;FAAGAAgAAAAhAI4h00HfAAAADAEAAA8AAAAAAAAAAAAAAAAA5wQAAGRycy9kb3ducmV2LnhtbFBL&#10;BQYAAAAABAAEAPMAAADzBQAAAAA=&#10;" fillcolor="#ffd869" strokecolor="#1f3763 [1604]" strokeweight="1pt"/></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p>
'@

$target = $d.Paragraphs(1).Range
$target.Collapse(1)
$target.InsertXML($newShapesXml)
